$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 159303
$ws.Range("C4").Value = 150349
$ws.Range("C5").Value = 8954
$ws.Range("C8").Value = 64.25
